$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, [string]$Text) {
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '69.050.14'
Set-TextValue $ws.Range('E2') '  -0.27%  '

Set-TextValue $ws.Range('D3') '3.805.77'
Set-TextValue $ws.Range('E3') '  +1.72%  '

Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.01%  '

Set-TextValue $ws.Range('D5') '600.81'
Set-TextValue $ws.Range('E5') '  -0.28%  '

Set-TextValue $ws.Range('D6') '163.90'
Set-TextValue $ws.Range('E6') '  -2.54%  '

Set-TextValue $ws.Range('D7') '3.802.77'
Set-TextValue $ws.Range('E7') '  +1.69%  '

Set-TextValue $ws.Range('E8') '  -0.03%  '

Set-TextValue $ws.Range('E9') '  -0.44%  '

Set-TextValue $ws.Range('E10') '  +1.44%  '

Set-TextValue $ws.Range('D11') '6.30'
Set-TextValue $ws.Range('E11') '  -1.88%  '

Set-TextValue $ws.Range('E12') '  -0.30%  '

Set-TextValue $ws.Range('D13') '37.08'
Set-TextValue $ws.Range('E13') '  -2.54%  '

Set-TextValue $ws.Range('E14') '  -0.67%  '

Set-TextValue $ws.Range('D15') '4.439.82'
Set-TextValue $ws.Range('E15') '  +1.71%  '

Set-TextValue $ws.Range('D16') '3.816.44'
Set-TextValue $ws.Range('E16') '  +2.06%  '

Set-TextValue $ws.Range('D17') '69.174.56'
Set-TextValue $ws.Range('E17') '  -0.04%  '

Set-TextValue $ws.Range('D18') '7.48'
Set-TextValue $ws.Range('E18') '  +2.37%  '

Set-TextValue $ws.Range('E19') '  +0.37%  '

Set-TextValue $ws.Range('D20') '11.43'
Set-TextValue $ws.Range('E20') '  +5.04%  '

Set-TextValue $ws.Range('D21') '17.22'
Set-TextValue $ws.Range('E21') '  +1.28%  '

Set-TextValue $ws.Range('D22') '487.12'
Set-TextValue $ws.Range('E22') '  -1.58%  '

Set-TextValue $ws.Range('D23') '0.719'
Set-TextValue $ws.Range('E23') '  -0.84%  '

Set-TextValue $ws.Range('D24') '0.0000160'
Set-TextValue $ws.Range('E24') '  +6.66%  '

Set-TextValue $ws.Range('D25') '84.62'
Set-TextValue $ws.Range('E25') '  -0.30%  '

Set-TextValue $ws.Range('E26') '  -2.55%  '

Set-TextValue $ws.Range('D27') '12.18'
Set-TextValue $ws.Range('E27') '  -0.85%  '

Set-TextValue $ws.Range('E28') '  -1.12%  '

Set-TextValue $ws.Range('E29') '  -0.09%  '

Set-TextValue $ws.Range('E30') '  -0.69%  '

Set-TextValue $ws.Range('D31') '8.00'
Set-TextValue $ws.Range('E31') '  -1.33%  '

Set-TextValue $ws.Range('E32') '  -4.91%  '

Set-TextValue $ws.Range('D33') '3.962.00'
Set-TextValue $ws.Range('E33') '  +1.93%  '

Set-TextValue $ws.Range('D34') '31.80'
Set-TextValue $ws.Range('E34') '  +0.56%  '

Set-TextValue $ws.Range('D35') '3.749.78'
Set-TextValue $ws.Range('E35') '  +2.07%  '

Set-TextValue $ws.Range('E36') '  -1.59%  '

Set-TextValue $ws.Range('E37') '  +0.73%  '

Set-TextValue $ws.Range('E38') '  +4.72%  '

Set-TextValue $ws.Range('D39') '5.86'
Set-TextValue $ws.Range('E39') '  -0.21%  '

Set-TextValue $ws.Range('D40') '1.00'
Set-TextValue $ws.Range('E40') '  +0.04%  '

Set-TextValue $ws.Range('E41') '  +1.49%  '

Set-TextValue $ws.Range('D42') '0.318'
Set-TextValue $ws.Range('E42') '  -1.73%  '

Set-TextValue $ws.Range('D43') '437.82'
Set-TextValue $ws.Range('E43') '  +0.78%  '

Set-TextValue $ws.Range('D44') '48.57'
Set-TextValue $ws.Range('E44') '  -0.21%  '

Set-TextValue $ws.Range('E45') '  -0.40%  '

Set-TextValue $ws.Range('E46') '  -0.01%  '

Set-TextValue $ws.Range('D47') '8.36'
Set-TextValue $ws.Range('E47') '  -1.32%  '

Set-TextValue $ws.Range('D48') '2.826.56'
Set-TextValue $ws.Range('E48') '  +1.85%  '

Set-TextValue $ws.Range('D49') '142.09'
Set-TextValue $ws.Range('E49') '  +0.88%  '

Set-TextValue $ws.Range('D50') '39.26'
Set-TextValue $ws.Range('E50') '  -2.76%  '

Set-TextValue $ws.Range('D51') '0.0351'
Set-TextValue $ws.Range('E51') '  -0.62%  '
